$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")
$ws.Activate()

# Formula used throughout column F (same pattern as all existing rows)
function Set-JsonFormula($row) {
    $f = '=_xlfn.CONCAT( ,A' + $row + ',": { ""worldId"": ",C' + $row + ',", ""name"": """,D' + $row + ',""", ""display"": """,E' + $row + ',""", ""areaId"": ",B' + $row + ',", },")'
    $ws.Range("F$row").Formula = $f
}

# --- Fill in the new Hollow Bastion RP locations for rows 70-76 ---

# Row 70 (id 69): HollowBastionLibrary3
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = "0x8"
$ws.Range("C70").Value = "0x5"
$ws.Range("D70").Value = "HollowBastionLibrary3"
$ws.Range("E70").Value = "Library 3"
Set-JsonFormula 70

# Row 71 (id 70): HollowBastionLibrary4
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "0x9"
$ws.Range("C71").Value = "0x5"
$ws.Range("D71").Value = "HollowBastionLibrary4"
$ws.Range("E71").Value = "Library 4"
Set-JsonFormula 71

# Row 72 (id 71): HollowBastionLibrary2
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "0x7"
$ws.Range("C72").Value = "0x5"
$ws.Range("D72").Value = "HollowBastionLibrary2"
$ws.Range("E72").Value = "Library 2"
Set-JsonFormula 72

# Row 73 (id 72): HollowBastionEntranceUpper
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "0x5"
$ws.Range("C73").Value = "0x5"
$ws.Range("D73").Value = "HollowBastionEntranceUpper"
$ws.Range("E73").Value = "Entrance Hall (Upper Level)"
Set-JsonFormula 73

# Row 74 (id 73): HollowBastionCrest
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "0xa"
$ws.Range("C74").Value = "0x5"
$ws.Range("D74").Value = "HollowBastionCrest"
$ws.Range("E74").Value = "Great Crest"
Set-JsonFormula 74

# Row 75 (id 74): HollowBastionTower
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "0xb"
$ws.Range("C75").Value = "0x5"
$ws.Range("D75").Value = "HollowBastionTower"
$ws.Range("E75").Value = "High Tower"
Set-JsonFormula 75

# Row 76 (id 75): HollowBastionGrandHall
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = "0xd"
$ws.Range("C76").Value = "0x5"
$ws.Range("D76").Value = "HollowBastionGrandHall"
$ws.Range("E76").Value = "Grand Hall(?)"
Set-JsonFormula 76

# --- New blank placeholder rows 77-85 (index-only + JSON formula, like sector rewards) ---
for ($row = 77; $row -le 85; $row++) {
    $ws.Range("A$row").Value = $row - 1
    Set-JsonFormula $row
}

# --- Restore the view state: scroll position + selected cell ---
$ws.Range("B77").Select()
$excel.ActiveWindow.ScrollRow = 46
